$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the cells in the same order the author typed them (this controls
# the order new entries land in the shared-strings table): columns A/B for
# both rows, then C/D, then E, then F; G (Status) last for each row.
$ws.Range("A11").Value = "TC010"
$ws.Range("B11").Value = "Valid and Broken Images"
$ws.Range("B12").Value = "Valid and Broken Links"
$ws.Range("A12").Value = "TC011"

$ws.Range("C11").Value = "Verify that both broken and non-broken images are correctly identified on the page."
$ws.Range("D11").Value = "Navigate to Broken Links-Images page."
$ws.Range("C12").Value = "Verify that the valid link redirects to the correct page and the broken link results in an error page."
$ws.Range("D12").Value = "Navigate to Broken Links-Images page."

$ws.Range("E11").Value = "1. Check whether the  'naturalWidth' Attribute is equal to 0 or not."
$ws.Range("E12").Value = "1. Open the links with Ctrl+Click.`n2. Verify the response code."

$ws.Range("F11").Value = "The valid image should have a naturalWidth greater than 0, meaning it is displayed correctly.`nThe broken image should have a naturalWidth of 0, indicating it is not displayed."
$ws.Range("F12").Value = "The valid link should open in a new tab and load successfully with an HTTP 200 status code.`nThe broken link should open in a new tab and return an HTTP 500 status code (indicating an error page)."

$ws.Range("G11").Value = "Pass"
$ws.Range("G12").Value = "Pass"

# Row heights for the new rows
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 60

# Wrap text for E and F columns on new rows (matches existing style of column E/F)
$ws.Range("E11:F12").WrapText = $true

# Center alignment for the Status (G) column on new rows
$ws.Range("G11:G12").HorizontalAlignment = -4108
$ws.Range("G11:G12").VerticalAlignment = -4108

# Columns C and D are "best fit" (bestFit/customWidth) - they grow to match
# the new (longer) content, mirroring Excel's automatic best-fit resizing.
$ws.Columns.Item(3).ColumnWidth = 88.3
$ws.Columns.Item(4).ColumnWidth = 35.0

# Update selection / view state
$null = $ws.Range("C6").Select()
